$d = $word.ActiveDocument

# Find the paragraph that contains the "LOQ4233..." requirement text,
# then remove the three paragraphs that follow it:
#   1. an empty paragraph
#   2. "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3. "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages.
#       Original theme under Creative Commons Attribution"
# leaving the remaining empty paragraph (and the page-break paragraph after it)
# untouched.

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r`a") -eq "LOQ4233: Gestão de Negócios (Requisito fraco)") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $startPara = $target.Next()
    $endPara = $startPara.Next().Next()
    $start = $startPara.Range.Start
    $end = $endPara.Range.End
    $r = $d.Range($start, $end)
    $r.Delete()
}
